# =====================================================================
# Generate Report for Handoff
#
# A new localization entry, 4f2cba96-e7c6-468c-b728-2ba6803ecb7d, is
# ready for handoff. It is inserted as row 8 on every worksheet (just
# before the existing d272bf14 entry), which pushes d272bf14 down to
# row 9 unchanged. Done via COM on all three sheets: Overview, zh-cn, de-de.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ===================== Sheet "Overview" =====================
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Rows.Item(8).Insert()

# New row 8: the 4f2cba96 entry that is ready for handoff
$ws1.Range('A8').Value2 = '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md'
$ws1.Range('B8').Value2 = 'Ready for handoff'
$ws1.Range('C8').Value2 = 'Ready for handoff'
$ws1.Range('D8').Value2 = '2016-46-18 02:46:36'

# Hyperlinks do not auto-follow the row-insert shift in this engine,
# so rebuild the whole collection in final, row-correct form.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/75c975e5f6c16b4199460c40b4a6d4062e7ba45d/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.md')
$ws1.Hyperlinks.Add($ws1.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/90d174cbaa3545334135e62dd5f473cd94d9b74e/e2e/0f28a0db-adcd-4868-8423-4140fde232f3.md', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.md')
$ws1.Hyperlinks.Add($ws1.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/222875f0-907e-479a-99bf-a8b94830f467.md', '', '', '222875f0-907e-479a-99bf-a8b94830f467.md')
$ws1.Hyperlinks.Add($ws1.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/6de7938d-9893-47d4-a56f-31dc0eac1cfe.md', '', '', '6de7938d-9893-47d4-a56f-31dc0eac1cfe.md')
$ws1.Hyperlinks.Add($ws1.Range('A6'), 'https://github.com/OpenLocalizationTest/oltest/blob/42162df6702f243acae5a83a2d76dfec92a7119b/e2e/9108f6ff-b6e6-4f65-9bec-cc42006e03af.md', '', '', '9108f6ff-b6e6-4f65-9bec-cc42006e03af.md')
$ws1.Hyperlinks.Add($ws1.Range('A7'), 'https://github.com/OpenLocalizationTest/oltest/blob/3f22aaa661a73359a2c16809f8a7f56406bb5015/e2e/44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md', '', '', '44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md')
$ws1.Hyperlinks.Add($ws1.Range('A8'), 'https://github.com/OpenLocalizationTest/oltest/blob/8a1c1e2d3f4b5a6c7d8e9f0a1b2c3d4e5f6a7b8c/e2e/4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md', '', '', '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md')
$ws1.Hyperlinks.Add($ws1.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/7d573aa996b0c8647155edbc2cc9874b57274480/e2e/d272bf14-eed9-4063-bdd0-417499bd8e8c.md', '', '', 'd272bf14-eed9-4063-bdd0-417499bd8e8c.md')

# ===================== Sheet "zh-cn" =====================
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(8).Insert()

# New row 8: the 4f2cba96 entry that is ready for handoff
$ws2.Range('A8').Value2 = '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md'
$ws2.Range('B8').Value2 = '.md'
$ws2.Range('C8').Value2 = 'Ready for handoff'
$ws2.Range('D8').Value2 = '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.9147ae4988c1bede7236c8b1eb19cafb80c2faf6.zh-cn.xlf'
$ws2.Range('E8').Value2 = '2016-03-18 02:46:28'
$ws2.Range('H8').Value2 = '0001-01-01 00:00:00'
$ws2.Range('I8').Value2 = 'Include'

# Hyperlinks do not auto-follow the row-insert shift in this engine,
# so rebuild the whole collection in final, row-correct form.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/75c975e5f6c16b4199460c40b4a6d4062e7ba45d/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.md')
$ws2.Hyperlinks.Add($ws2.Range('B2'), 'https://github.com/OpenLocalizationTest/oltest/blob/75c975e5f6c16b4199460c40b4a6d4062e7ba45d/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/459f3ffb4f8221b8359894be1f71da584bca470a/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.zh-cn.xlf', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/27d1f17ee713583056794eac4aa422ec52d851d0/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.md')
$ws2.Hyperlinks.Add($ws2.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/c01d3619e399544b15209536ef7c673513f104e0/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.zh-cn.xlf', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/90d174cbaa3545334135e62dd5f473cd94d9b74e/e2e/0f28a0db-adcd-4868-8423-4140fde232f3.md', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.md')
$ws2.Hyperlinks.Add($ws2.Range('B3'), 'https://github.com/OpenLocalizationTest/oltest/blob/90d174cbaa3545334135e62dd5f473cd94d9b74e/e2e/0f28a0db-adcd-4868-8423-4140fde232f3.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4794ad6838b446420731a61f6433e155335b1fbd/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.zh-cn.xlf', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/d04cd1b48fdec7cd8cc306d629b5f1b491bbacae/e2e/0f28a0db-adcd-4868-8423-4140fde232f3.md', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.md')
$ws2.Hyperlinks.Add($ws2.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/87a6b6b104c373c615e12c925d582353f80bbaea/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.zh-cn.xlf', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/222875f0-907e-479a-99bf-a8b94830f467.md', '', '', '222875f0-907e-479a-99bf-a8b94830f467.md')
$ws2.Hyperlinks.Add($ws2.Range('B4'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/222875f0-907e-479a-99bf-a8b94830f467.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b58d48b4b51d881d18df2f827562167da1273289/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/222875f0-907e-479a-99bf-a8b94830f467.ecfdadfe873a6c7dbbd2d395a9ab3a30be6ff04f.zh-cn.xlf', '', '', '222875f0-907e-479a-99bf-a8b94830f467.ecfdadfe873a6c7dbbd2d395a9ab3a30be6ff04f.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/6de7938d-9893-47d4-a56f-31dc0eac1cfe.md', '', '', '6de7938d-9893-47d4-a56f-31dc0eac1cfe.md')
$ws2.Hyperlinks.Add($ws2.Range('B5'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/6de7938d-9893-47d4-a56f-31dc0eac1cfe.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D5'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b58d48b4b51d881d18df2f827562167da1273289/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/6de7938d-9893-47d4-a56f-31dc0eac1cfe.ee928456acfa01a535c9251599edbeaebf3cfc6c.zh-cn.xlf', '', '', '6de7938d-9893-47d4-a56f-31dc0eac1cfe.ee928456acfa01a535c9251599edbeaebf3cfc6c.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A6'), 'https://github.com/OpenLocalizationTest/oltest/blob/42162df6702f243acae5a83a2d76dfec92a7119b/e2e/9108f6ff-b6e6-4f65-9bec-cc42006e03af.md', '', '', '9108f6ff-b6e6-4f65-9bec-cc42006e03af.md')
$ws2.Hyperlinks.Add($ws2.Range('B6'), 'https://github.com/OpenLocalizationTest/oltest/blob/42162df6702f243acae5a83a2d76dfec92a7119b/e2e/9108f6ff-b6e6-4f65-9bec-cc42006e03af.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D6'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e4c59540937d4b0995a495c6832c637da886bde/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/9108f6ff-b6e6-4f65-9bec-cc42006e03af.f8b59048bc8c7459296b55851d3372f6f026be07.zh-cn.xlf', '', '', '9108f6ff-b6e6-4f65-9bec-cc42006e03af.f8b59048bc8c7459296b55851d3372f6f026be07.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A7'), 'https://github.com/OpenLocalizationTest/oltest/blob/3f22aaa661a73359a2c16809f8a7f56406bb5015/e2e/44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md', '', '', '44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md')
$ws2.Hyperlinks.Add($ws2.Range('B7'), 'https://github.com/OpenLocalizationTest/oltest/blob/3f22aaa661a73359a2c16809f8a7f56406bb5015/e2e/44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D7'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa46be93d0af565759715f5cec1a2218753e18c5/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/44e48f6d-14e0-46ea-9fd1-e0dacc693fab.3de409930bc3354428a9bd2a0523794cf539c963.zh-cn.xlf', '', '', '44e48f6d-14e0-46ea-9fd1-e0dacc693fab.3de409930bc3354428a9bd2a0523794cf539c963.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A8'), 'https://github.com/OpenLocalizationTest/oltest/blob/8a1c1e2d3f4b5a6c7d8e9f0a1b2c3d4e5f6a7b8c/e2e/4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md', '', '', '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md')
$ws2.Hyperlinks.Add($ws2.Range('B8'), 'https://github.com/OpenLocalizationTest/oltest/blob/8a1c1e2d3f4b5a6c7d8e9f0a1b2c3d4e5f6a7b8c/e2e/4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D8'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a1c1e2d3f4b5a6c7d8e9f0a1b2c3d4e5f6a7b8c/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/4f2cba96-e7c6-468c-b728-2ba6803ecb7d.9147ae4988c1bede7236c8b1eb19cafb80c2faf6.zh-cn.xlf', '', '', '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.9147ae4988c1bede7236c8b1eb19cafb80c2faf6.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/7d573aa996b0c8647155edbc2cc9874b57274480/e2e/d272bf14-eed9-4063-bdd0-417499bd8e8c.md', '', '', 'd272bf14-eed9-4063-bdd0-417499bd8e8c.md')
$ws2.Hyperlinks.Add($ws2.Range('B9'), 'https://github.com/OpenLocalizationTest/oltest/blob/7d573aa996b0c8647155edbc2cc9874b57274480/e2e/d272bf14-eed9-4063-bdd0-417499bd8e8c.md', '', '', '.md')
$ws2.Hyperlinks.Add($ws2.Range('D9'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/17bdd99539566ff19e359028fd4b571054a9c55c/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d272bf14-eed9-4063-bdd0-417499bd8e8c.2cbf6f293531aae3538c7e4c40b16a3334bbc068.zh-cn.xlf', '', '', 'd272bf14-eed9-4063-bdd0-417499bd8e8c.2cbf6f293531aae3538c7e4c40b16a3334bbc068.zh-cn.xlf')

# ===================== Sheet "de-de" =====================
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(8).Insert()

# New row 8: the 4f2cba96 entry that is ready for handoff
$ws3.Range('A8').Value2 = '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md'
$ws3.Range('B8').Value2 = '.md'
$ws3.Range('C8').Value2 = 'Ready for handoff'
$ws3.Range('D8').Value2 = '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.9147ae4988c1bede7236c8b1eb19cafb80c2faf6.de-de.xlf'
$ws3.Range('E8').Value2 = '2016-03-18 02:46:36'
$ws3.Range('H8').Value2 = '0001-01-01 00:00:00'
$ws3.Range('I8').Value2 = 'Include'

# Hyperlinks do not auto-follow the row-insert shift in this engine,
# so rebuild the whole collection in final, row-correct form.
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/75c975e5f6c16b4199460c40b4a6d4062e7ba45d/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.md')
$ws3.Hyperlinks.Add($ws3.Range('B2'), 'https://github.com/OpenLocalizationTest/oltest/blob/75c975e5f6c16b4199460c40b4a6d4062e7ba45d/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d6e0b46cf74a57af7a80e8225267536a9b9cb93/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.de-de.xlf', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/d4e4818fc48f4a7f6067c47f1731c4b3c6cb9878/e2e/dd233af3-56c7-4c64-a07e-cf855ca0f82d.md', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.md')
$ws3.Hyperlinks.Add($ws3.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/7a16cff070bb5e8158d2cdf598bb97a3e3f0f852/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.de-de.xlf', '', '', 'dd233af3-56c7-4c64-a07e-cf855ca0f82d.a3e3d5d1a309f9944f6103f007df5fd5c012a303.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/90d174cbaa3545334135e62dd5f473cd94d9b74e/e2e/0f28a0db-adcd-4868-8423-4140fde232f3.md', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.md')
$ws3.Hyperlinks.Add($ws3.Range('B3'), 'https://github.com/OpenLocalizationTest/oltest/blob/90d174cbaa3545334135e62dd5f473cd94d9b74e/e2e/0f28a0db-adcd-4868-8423-4140fde232f3.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34475ed74c232803a99a3f9315ddb0f37be6b8bf/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.de-de.xlf', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/052f5cee4261568e293aaaacc39ad73381685f18/e2e/0f28a0db-adcd-4868-8423-4140fde232f3.md', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.md')
$ws3.Hyperlinks.Add($ws3.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/9e0be39e1b006afab85b786a926c8f1278dcaa5e/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.de-de.xlf', '', '', '0f28a0db-adcd-4868-8423-4140fde232f3.fead972be1f9183c09474f536144deb0107dbca3.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/222875f0-907e-479a-99bf-a8b94830f467.md', '', '', '222875f0-907e-479a-99bf-a8b94830f467.md')
$ws3.Hyperlinks.Add($ws3.Range('B4'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/222875f0-907e-479a-99bf-a8b94830f467.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fe2164f78855c6857b6d526e280e05b5f570e03/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/222875f0-907e-479a-99bf-a8b94830f467.ecfdadfe873a6c7dbbd2d395a9ab3a30be6ff04f.de-de.xlf', '', '', '222875f0-907e-479a-99bf-a8b94830f467.ecfdadfe873a6c7dbbd2d395a9ab3a30be6ff04f.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/6de7938d-9893-47d4-a56f-31dc0eac1cfe.md', '', '', '6de7938d-9893-47d4-a56f-31dc0eac1cfe.md')
$ws3.Hyperlinks.Add($ws3.Range('B5'), 'https://github.com/OpenLocalizationTest/oltest/blob/5c1492444e31b7282a614034279cfce99bc22ba8/e2e/6de7938d-9893-47d4-a56f-31dc0eac1cfe.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D5'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fe2164f78855c6857b6d526e280e05b5f570e03/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/6de7938d-9893-47d4-a56f-31dc0eac1cfe.ee928456acfa01a535c9251599edbeaebf3cfc6c.de-de.xlf', '', '', '6de7938d-9893-47d4-a56f-31dc0eac1cfe.ee928456acfa01a535c9251599edbeaebf3cfc6c.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A6'), 'https://github.com/OpenLocalizationTest/oltest/blob/42162df6702f243acae5a83a2d76dfec92a7119b/e2e/9108f6ff-b6e6-4f65-9bec-cc42006e03af.md', '', '', '9108f6ff-b6e6-4f65-9bec-cc42006e03af.md')
$ws3.Hyperlinks.Add($ws3.Range('B6'), 'https://github.com/OpenLocalizationTest/oltest/blob/42162df6702f243acae5a83a2d76dfec92a7119b/e2e/9108f6ff-b6e6-4f65-9bec-cc42006e03af.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D6'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/086ad46af769e9055a5c3f7664ee8117073edc32/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/9108f6ff-b6e6-4f65-9bec-cc42006e03af.f8b59048bc8c7459296b55851d3372f6f026be07.de-de.xlf', '', '', '9108f6ff-b6e6-4f65-9bec-cc42006e03af.f8b59048bc8c7459296b55851d3372f6f026be07.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A7'), 'https://github.com/OpenLocalizationTest/oltest/blob/3f22aaa661a73359a2c16809f8a7f56406bb5015/e2e/44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md', '', '', '44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md')
$ws3.Hyperlinks.Add($ws3.Range('B7'), 'https://github.com/OpenLocalizationTest/oltest/blob/3f22aaa661a73359a2c16809f8a7f56406bb5015/e2e/44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D7'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6be515222e207659f5ab54c762c4bcdd28c09360/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/44e48f6d-14e0-46ea-9fd1-e0dacc693fab.3de409930bc3354428a9bd2a0523794cf539c963.de-de.xlf', '', '', '44e48f6d-14e0-46ea-9fd1-e0dacc693fab.3de409930bc3354428a9bd2a0523794cf539c963.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A8'), 'https://github.com/OpenLocalizationTest/oltest/blob/8a1c1e2d3f4b5a6c7d8e9f0a1b2c3d4e5f6a7b8c/e2e/4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md', '', '', '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md')
$ws3.Hyperlinks.Add($ws3.Range('B8'), 'https://github.com/OpenLocalizationTest/oltest/blob/8a1c1e2d3f4b5a6c7d8e9f0a1b2c3d4e5f6a7b8c/e2e/4f2cba96-e7c6-468c-b728-2ba6803ecb7d.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D8'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a1c1e2d3f4b5a6c7d8e9f0a1b2c3d4e5f6a7b8c/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/4f2cba96-e7c6-468c-b728-2ba6803ecb7d.9147ae4988c1bede7236c8b1eb19cafb80c2faf6.de-de.xlf', '', '', '4f2cba96-e7c6-468c-b728-2ba6803ecb7d.9147ae4988c1bede7236c8b1eb19cafb80c2faf6.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/7d573aa996b0c8647155edbc2cc9874b57274480/e2e/d272bf14-eed9-4063-bdd0-417499bd8e8c.md', '', '', 'd272bf14-eed9-4063-bdd0-417499bd8e8c.md')
$ws3.Hyperlinks.Add($ws3.Range('B9'), 'https://github.com/OpenLocalizationTest/oltest/blob/7d573aa996b0c8647155edbc2cc9874b57274480/e2e/d272bf14-eed9-4063-bdd0-417499bd8e8c.md', '', '', '.md')
$ws3.Hyperlinks.Add($ws3.Range('D9'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d806dbb1cc32907545157a1e2fe6cce100a99092/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d272bf14-eed9-4063-bdd0-417499bd8e8c.2cbf6f293531aae3538c7e4c40b16a3334bbc068.de-de.xlf', '', '', 'd272bf14-eed9-4063-bdd0-417499bd8e8c.2cbf6f293531aae3538c7e4c40b16a3334bbc068.de-de.xlf')

